# Mark the checklist items as passed.
# The "Чек-лист проверок" sheet has a two-column layout: column A holds the
# check description, column B holds its status ("Passed" / "Fail").
# Every check row below the section headers gets its status set to "Passed".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$passedRows = 4,5,6,7,9,10,11,12,13,14,15,16,17,18,19,20,21,22

foreach ($r in $passedRows) {
    $ws.Cells.Item($r, 2).Value = "Passed"
}

$wb.Save()
